$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of daily covid-19 case data for San Diego County, Mar 28, 2020 (row 24)
$rowNum = 24
$rowData = [ordered]@{
    "A"  = 43918
    "B"  = 488
    "C"  = 3
    "D"  = 0
    "E"  = 5
    "F"  = 1
    "G"  = 102
    "H"  = 5
    "I"  = 118
    "J"  = 14
    "K"  = 85
    "L"  = 15
    "M"  = 74
    "N"  = 21
    "O"  = 36
    "P"  = 12
    "Q"  = 38
    "R"  = 13
    "S"  = 25
    "T"  = 15
    "U"  = 2
    "V"  = 0
    "W"  = 206
    "X"  = 280
    "Y"  = 2
    "Z"  = 96
    "AA" = 42
    "AB" = 7
    "AC" = 19
    "AD" = 24
    "AE" = 1
    "AF" = 6
    "AG" = 26
    "AH" = 16
    "AI" = 9
    "AJ" = 3
    "AK" = 2
    "AL" = 8
    "AM" = 11
    "AN" = 7
    "AO" = 284
    "AP" = 3
    "AQ" = 2
    "AR" = 1
    "AS" = 10
    "AT" = 2
    "AU" = 4
    "AW" = 2
    "AX" = 4
    "AY" = 7
    "BA" = 9
}

foreach ($col in $rowData.Keys) {
    $cellRef = "{0}{1}" -f $col, $rowNum
    $ws.Range($cellRef).Value = $rowData[$col]
}

# Mirror the scrolled/selected cell left after entering the new row
$ws.Range("X24").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$win.ScrollRow = 1

$wb.Save()
